$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 63
$ws.Range("F5").Value = 227
$ws.Range("F6").Value = 264
$ws.Range("F7").Value = 62
$ws.Range("F13").Value = 2186
$ws.Range("F16").Value = 508
$ws.Range("F17").Value = 499
$ws.Range("F18").Value = 154
$ws.Range("F20").Value = 38
$ws.Range("F21").Value = 44
$ws.Range("F22").Value = 1674
$ws.Range("F23").Value = 3838
$ws.Range("F24").Value = 28
$ws.Range("F25").Value = 60
$ws.Range("F27").Value = 1147
$ws.Range("F28").Value = 212
$ws.Range("G28").Value = 46.6
$ws.Range("F29").Value = 2037
$ws.Range("F32").Value = 83
$ws.Range("F34").Value = 414
$ws.Range("F35").Value = 456
$ws.Range("F36").Value = 665
$ws.Range("F38").Value = 397

# --- Sheet "演出" (shows) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 23

# --- Sheet "全部类型" (all types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 63
$ws.Range("F5").Value = 227
$ws.Range("F6").Value = 264
$ws.Range("F7").Value = 62
$ws.Range("F13").Value = 2186
$ws.Range("F15").Value = 23
$ws.Range("F17").Value = 508
$ws.Range("F18").Value = 499
$ws.Range("F19").Value = 154
$ws.Range("F21").Value = 38
$ws.Range("F22").Value = 44
$ws.Range("F23").Value = 1674
$ws.Range("F24").Value = 3838
$ws.Range("F25").Value = 28
$ws.Range("F26").Value = 60
$ws.Range("F28").Value = 1147
$ws.Range("F29").Value = 212
$ws.Range("G29").Value = 46.6
$ws.Range("F30").Value = 2037
$ws.Range("F33").Value = 83
$ws.Range("F35").Value = 414
$ws.Range("F36").Value = 456
$ws.Range("F37").Value = 665
$ws.Range("F39").Value = 397
